$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "s"
$ws.Range("B1").Value = "a"
$ws.Range("C1").Value = "y"

$ws.Range("C1").Select()
